$wb = $excel.ActiveWorkbook

# Add the new "Tasks" sheet. A freshly-added sheet is inserted as the first
# sheet and becomes the active/tab-selected sheet, matching the target
# workbook order (Tasks, Tables, References).
$ws = $wb.Worksheets.Add()
$ws.Name = "Tasks"

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Description"

# Reuse the existing "thin box border" formatting (style already used on the
# References sheet) so no new border gets created, then layer the yellow
# header fill on top - this reproduces the single new cellXfs entry
# (fontId=0, fillId=2, borderId=1) from the diff instead of inventing a new
# border/fill pair.
$refs = $wb.Worksheets.Item("References")
$refs.Range("A2").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$ws.Range("A1:B1").Interior.Color = 65535

$ws.Columns.Item(2).ColumnWidth = 19.6

$null = $ws.Range("A2").Select()
